# Apply the DataFormat.xlsx restructuring described by the commit:
# "improve structure, sensor recv by FSM"
#
# - Row 1 / Row 2 header cells keep their look (border, centered) - no visual
#   change there, only the K2:N2 block loses a (already invisible) fill flag.
# - New explanatory notes are added in row 5, 6, 7 and a new row 8.
# - The L5 note text changes, and the selection moves from L9 to L6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- K2:N2: drop the (no-op) fill flag so the cells fall back to the plain
#     bordered style, matching the trimmed cellXfs table in the target file.
$ws.Range("K2:N2").Interior.Pattern = -4142   # xlPatternNone

# --- New note cells (row 5) ------------------------------------------------
$ws.Range("A5").Value = "所有数据位均按uint16_t解读"
$ws.Range("O5").Value = "单位摄氏度℃"

# --- New note cell (row 6) --------------------------------------------------
$ws.Range("L6").Value = "压力单位为Pa"

# --- New note cells (row 7) -------------------------------------------------
$ws.Range("L7").Value = "水深为乘1000后结果"
$ws.Range("O7").Value = "单位米m"

# --- New note cell (row 8, brand-new row) -----------------------------------
$ws.Range("L8").Value = "舱温舱湿为乘1000后结果"

# --- Existing L5 note text changes (set last so it lands on a fresh shared
#     string slot instead of reusing the one just freed above) --------------
$ws.Range("L5").Value = "水温为乘100后结果"

# --- Selection moves from L9 to L6 ------------------------------------------
$ws.Range("L6").Select()
